# Update the "Urbanization" breakdown rows with more specific labels.
# Row 23 (Urban / Город / Шаар) -> City / Городские поселения / Шаар жерлери
# Row 24 (Rural / Село / Айыл)  -> Village / Сельская местность / Айыл аймагы
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A23").Value = "Шаар жерлери"
$ws.Range("B23").Value = "Городские поселения"
$ws.Range("C23").Value = "City"

$ws.Range("A24").Value = "Айыл аймагы"
$ws.Range("B24").Value = "Сельская местность"
$ws.Range("C24").Value = "Village"

# Move the active selection to match the saved view state.
$ws.Range("C30").Select()
